$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Range("H99").Value = 3932.182
$ws.Range("I99").Value = 3556.4
$ws.Range("J99").Value = 4245.3335
$ws.Range("K99").Value = 10669.2
$ws.Range("L99").Value = 12736.0005
$ws.Range("M99").Value = -9171.200000000001
$ws.Range("N99").Value = -15732.0005
# Row 135
$ws.Range("H135").Value = 1968.6765
$ws.Range("I135").Value = 1064.5416
$ws.Range("J135").Value = 4138.6
$ws.Range("K135").Value = 9580.874400000001
$ws.Range("L135").Value = 37247.4
$ws.Range("M135").Value = -7045.874400000001
$ws.Range("N135").Value = -42317.4
# Row 138
$ws.Range("H138").Value = 3850.544
$ws.Range("I138").Value = 3705.111
$ws.Range("K138").Value = 11115.333
$ws.Range("M138").Value = -5975.332999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10994.435
$ws.Range("I32").Value = 8046.86
$ws.Range("J32").Value = 24995.416
$ws.Range("K32").Value = 8046.86
$ws.Range("L32").Value = 24995.416
$ws.Range("M32").Value = -7759.86
$ws.Range("N32").Value = -25569.416
# Row 61
$ws.Range("H61").Value = 5381.0415
$ws.Range("I61").Value = 4857.25
$ws.Range("K61").Value = 4857.25
$ws.Range("M61").Value = -4645.25
# Row 97
$ws.Range("H97").Value = 1123310.5
$ws.Range("J97").Value = 9592.362999999999
$ws.Range("L97").Value = 9592.362999999999
$ws.Range("N97").Value = -10584.363
# Row 132
$ws.Range("H132").Value = 25043.871
$ws.Range("I132").Value = 1701.3793
$ws.Range("J132").Value = 92737.10000000001
$ws.Range("K132").Value = 5104.1379
$ws.Range("L132").Value = 278211.3
$ws.Range("M132").Value = -2574.1379
$ws.Range("N132").Value = -283271.3
# Row 136
$ws.Range("H136").Value = 5381.0415
$ws.Range("I136").Value = 4857.25
$ws.Range("K136").Value = 14571.75
$ws.Range("M136").Value = -12021.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 12040969
$ws.Range("I86").Value = 19699708
$ws.Range("K86").Value = 19699708
$ws.Range("M86").Value = -19698585
# Row 89
$ws.Range("H89").Value = 12040969
$ws.Range("I89").Value = 19699708
$ws.Range("K89").Value = 98498540
$ws.Range("M89").Value = -98492924
# Row 94
$ws.Range("H94").Value = 3680477
$ws.Range("I94").Value = 5209991.5
$ws.Range("J94").Value = 9642.5
$ws.Range("K94").Value = 5209991.5
$ws.Range("L94").Value = 9642.5
$ws.Range("M94").Value = -5209540.5
$ws.Range("N94").Value = -10544.5
# Row 107
$ws.Range("H107").Value = 4468973.5
$ws.Range("I107").Value = 7145422.5
$ws.Range("K107").Value = 7145422.5
$ws.Range("M107").Value = -7143502.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 17831.895
$ws.Range("I31").Value = 2692.7908
$ws.Range("K31").Value = 2692.7908
$ws.Range("M31").Value = -2397.7908
# Row 34
$ws.Range("H34").Value = 17831.895
$ws.Range("I34").Value = 2692.7908
$ws.Range("K34").Value = 2692.7908
$ws.Range("M34").Value = -2490.7908
# Row 105
$ws.Range("H105").Value = 809.2727
$ws.Range("I105").Value = 786.55554
$ws.Range("K105").Value = 786.55554
$ws.Range("M105").Value = 960.44446
# Row 141
$ws.Range("H141").Value = 250714.92
$ws.Range("J141").Value = 250714.92
$ws.Range("L141").Value = 250714.92
$ws.Range("N141").Value = -261074.92

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 79143.92
$ws.Range("I5").Value = 1159
$ws.Range("J5").Value = 203919.8
$ws.Range("K5").Value = 3477
$ws.Range("L5").Value = 611759.3999999999
$ws.Range("M5").Value = -3365
$ws.Range("N5").Value = -611983.3999999999
# Row 131
$ws.Range("H131").Value = 18945140
$ws.Range("I131").Value = 16667491
$ws.Range("J131").Value = 19615036
$ws.Range("K131").Value = 50002473
$ws.Range("L131").Value = 58845108
$ws.Range("M131").Value = -49997433
$ws.Range("N131").Value = -58855188
# Row 135
$ws.Range("H135").Value = 79143.92
$ws.Range("I135").Value = 1159
$ws.Range("J135").Value = 203919.8
$ws.Range("K135").Value = 10431
$ws.Range("L135").Value = 1835278.2
$ws.Range("M135").Value = -7896
$ws.Range("N135").Value = -1840348.2
# Row 137
$ws.Range("H137").Value = 6533.8335
$ws.Range("I137").Value = 3526.375
$ws.Range("J137").Value = 8939.799999999999
$ws.Range("K137").Value = 10579.125
$ws.Range("L137").Value = 26819.4
$ws.Range("M137").Value = -5479.125
$ws.Range("N137").Value = -37019.39999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1204292.8
$ws.Range("I80").Value = 1481749
$ws.Range("J80").Value = 418166.5
$ws.Range("K80").Value = 1481749
$ws.Range("L80").Value = 418166.5
$ws.Range("M80").Value = -1480751
$ws.Range("N80").Value = -420162.5
# Row 83
$ws.Range("H83").Value = 1204292.8
$ws.Range("I83").Value = 1481749
$ws.Range("J83").Value = 418166.5
$ws.Range("K83").Value = 7408745
$ws.Range("L83").Value = 2090832.5
$ws.Range("M83").Value = -7403753
$ws.Range("N83").Value = -2100816.5
# Row 102
$ws.Range("H102").Value = 4218222.5
$ws.Range("I102").Value = 5849153.5
$ws.Range("K102").Value = 5849153.5
$ws.Range("M102").Value = -5847531.5
# Row 132
$ws.Range("H132").Value = 3098.2273
$ws.Range("I132").Value = 2737.2727
$ws.Range("J132").Value = 4903
$ws.Range("K132").Value = 8211.8181
$ws.Range("L132").Value = 14709
$ws.Range("M132").Value = -5681.8181
$ws.Range("N132").Value = -19769

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 44845.555
$ws.Range("I10").Value = 175.08333
$ws.Range("J10").Value = 134186.5
$ws.Range("K10").Value = 175.08333
$ws.Range("L10").Value = 134186.5
$ws.Range("M10").Value = -35.08332999999999
$ws.Range("N10").Value = -134466.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8127.9062
$ws.Range("I62").Value = 4919.2
$ws.Range("K62").Value = 4919.2
$ws.Range("M62").Value = -4295.2
# Row 65
$ws.Range("H65").Value = 8127.9062
$ws.Range("I65").Value = 4919.2
$ws.Range("K65").Value = 24596
$ws.Range("M65").Value = -21476
# Row 122
$ws.Range("H122").Value = 3799.9333
$ws.Range("I122").Value = 2174.4783
$ws.Range("K122").Value = 6523.4349
$ws.Range("M122").Value = -4073.4349
# Row 132
$ws.Range("H132").Value = 15091205
$ws.Range("I132").Value = 18185292
$ws.Range("J132").Value = 909977.7
$ws.Range("K132").Value = 54555876
$ws.Range("L132").Value = 2729933.1
$ws.Range("M132").Value = -54553346
$ws.Range("N132").Value = -2734993.1
# Row 139
$ws.Range("H139").Value = 68735.25
$ws.Range("J139").Value = 68735.25
$ws.Range("L139").Value = 68735.25
$ws.Range("N139").Value = -79015.25
